## edit.ps1
## Adds 19 new coordinate / crime-incident rows (127-145) to Sheet1,
## matching the "Create coordinate data" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stamp the new rows with the same look & feel (fonts/number formats) as
#        the existing data rows, by cloning the formatting of the last data row.
$ws.Range("A126:E126").Copy()
$ws.Range("A127:E145").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Data rows use a taller row height (20) than the sheet default (17)
$ws.Range("A127:A145").EntireRow.RowHeight = 20

# --- 2. Fill in the coordinate / radius / type / grade values for each new row

$ws.Range("A127").Value = 36.350878967494801
$ws.Range("B127").Value = 127.419032028638
$ws.Range("C127").Value = 25
$ws.Range("D127").Value = "Violence"
$ws.Range("E127").Value = 3

$ws.Range("A128").Value = 36.3544926646439
$ws.Range("B128").Value = 127.411252747918
$ws.Range("C128").Value = 85
$ws.Range("D128").Value = "Violence"
$ws.Range("E128").Value = 3

$ws.Range("A129").Value = 36.354428204745297
$ws.Range("B129").Value = 127.41426045663501
$ws.Range("C129").Value = 60
$ws.Range("D129").Value = "Violence"
$ws.Range("E129").Value = 3

$ws.Range("A130").Value = 36.357227060146798
$ws.Range("B130").Value = 127.407501366788
$ws.Range("C130").Value = 30
$ws.Range("D130").Value = "Violence"
$ws.Range("E130").Value = 3

$ws.Range("A131").Value = 36.350678444258598
$ws.Range("B131").Value = 127.398621768656
$ws.Range("C131").Value = 150
$ws.Range("D131").Value = "Violence"
$ws.Range("E131").Value = 1

$ws.Range("A132").Value = 36.347932886949202
$ws.Range("B132").Value = 127.400412452678
$ws.Range("C132").Value = 25
$ws.Range("D132").Value = "Violence"
$ws.Range("E132").Value = 2

$ws.Range("A133").Value = 36.346291357852998
$ws.Range("B133").Value = 127.400827358941
$ws.Range("C133").Value = 30
$ws.Range("D133").Value = "Violence"
$ws.Range("E133").Value = 2

$ws.Range("A134").Value = 36.341987395811799
$ws.Range("B134").Value = 127.39973594383601
$ws.Range("C134").Value = 35
$ws.Range("D134").Value = "Violence"
$ws.Range("E134").Value = 2

$ws.Range("A135").Value = 36.340668195068503
$ws.Range("B135").Value = 127.400776264754
$ws.Range("C135").Value = 25
$ws.Range("D135").Value = "Violence"
$ws.Range("E135").Value = 2

$ws.Range("A136").Value = 36.338643501367301
# Longitude kept as text (trailing space) just like in the source data
$ws.Range("B136").Value = "'127.40789464214636 "
$ws.Range("A126").Copy()
$ws.Range("B136").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("C136").Value = 55
$ws.Range("D136").Value = "Violence"
$ws.Range("E136").Value = 2

$ws.Range("A137").Value = 36.336652179017499
$ws.Range("B137").Value = 127.41044608857899
$ws.Range("C137").Value = 100
$ws.Range("D137").Value = "Violence"
$ws.Range("E137").Value = 1

$ws.Range("A138").Value = 36.335149505512703
$ws.Range("B138").Value = 127.41239852353699
$ws.Range("C138").Value = 100
$ws.Range("D138").Value = "Violence"
$ws.Range("E138").Value = 1

$ws.Range("A139").Value = 36.335964157460602
$ws.Range("B139").Value = 127.411355818646
$ws.Range("C139").Value = 100
$ws.Range("D139").Value = "Violence"
$ws.Range("E139").Value = 1

$ws.Range("A140").Value = 36.336830247570603
$ws.Range("B140").Value = 127.416283550036
$ws.Range("C140").Value = 60
$ws.Range("D140").Value = "Violence"
$ws.Range("E140").Value = 2

$ws.Range("A141").Value = 36.338552853065302
$ws.Range("B141").Value = 127.408050110706
$ws.Range("C141").Value = 75
$ws.Range("D141").Value = "Violence"
$ws.Range("E141").Value = 1

$ws.Range("A142").Value = 36.3301994324584
$ws.Range("B142").Value = 127.420903732971
$ws.Range("C142").Value = 140
$ws.Range("D142").Value = "Violence"
$ws.Range("E142").Value = 1

$ws.Range("A143").Value = 36.338208014667998
$ws.Range("B143").Value = 127.429311908264
$ws.Range("C143").Value = 85
$ws.Range("D143").Value = "Violence"
$ws.Range("E143").Value = 1

$ws.Range("A144").Value = 36.337449240890102
$ws.Range("B144").Value = 127.429808976413
$ws.Range("C144").Value = 35
$ws.Range("D144").Value = "Violence"
$ws.Range("E144").Value = 1

$ws.Range("A145").Value = 36.344495654171901
$ws.Range("B145").Value = 127.44236848637399
$ws.Range("C145").Value = 60
$ws.Range("D145").Value = "Violence"
$ws.Range("E145").Value = 1

# --- 3. Update the view so the newly added rows are in frame, mirroring the
#        scrolled/selected state captured in the saved workbook.
$ws.Range("F128").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 112 } catch {}

